$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 0
$ws1.Range("F2").Value = 24
$ws1.Range("G2").Value = 100
$ws1.Range("H2").Value = 8.800000000000001

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E2").Value = 24

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 0
$ws3.Range("F2").Value = 24
$ws3.Range("G2").Value = 100
$ws3.Range("H2").Value = 8.800000000000001
